$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered/top-aligned) from H1 to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add the new data values
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 10
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
